$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the sync diff: cell address -> new value
$updates = [ordered]@{
    'D2' = '29.969.08'
    'E2' = '  -0.42%  '
    'D3' = '1.876.08'
    'D4' = '0.9996'
    'E4' = '  -0.08%  '
    'D5' = '242.27'
    'E5' = '  -3.57%  '
    'D6' = '0.9995'
    'E6' = '  -0.09%  '
    'D7' = '0.4926'
    'E7' = '  -3.66%  '
    'D8' = '0.2918'
    'E8' = '  -1.71%  '
    'D9' = '0.06632'
    'E9' = '  -2.75%  '
    'D10' = '1.877.98'
    'E10' = '  -1.66%  '
    'D11' = '16.72'
    'E11' = '  -3.04%  '
    'D12' = '0.07245'
    'E12' = '  -1.62%  '
    'D13' = '0.6654'
    'E13' = '  -3.86%  '
    'D14' = '86.20'
    'E14' = '  -0.57%  '
    'D15' = '4.886'
    'E15' = '  +0.26%  '
    'D16' = '29.937.61'
    'E16' = '  -0.54%  '
    'D17' = '0.000007863'
    'E17' = '  -4.38%  '
    'D18' = '0.9998'
    'E18' = '  -0.01%  '
    'D19' = '12.75'
    'E19' = '  -1.83%  '
    'D20' = '2.120.70'
    'E20' = '  -1.69%  '
    'D21' = '0.9992'
    'E21' = '  -0.09%  '
    'D22' = '4.769'
    'E22' = '  -0.93%  '
    'D23' = '5.763'
    'E23' = '  +0.91%  '
    'D24' = '9.043'
    'E24' = '  -1.53%  '
    'D25' = '149.72'
    'E25' = '  +2.00%  '
    'D26' = '141.20'
    'E26' = '  +4.45%  '
    'D27' = '17.00'
    'E27' = '  -0.10%  '
    'D28' = '1.912'
    'E28' = '  -4.35%  '
    'D29' = '1.392'
    'E29' = '  +0.07%  '
    'D30' = '4.186'
    'E30' = '  -1.05%  '
    'D31' = '0.08732'
    'E31' = '  -1.02%  '
    'D32' = '3.959'
    'E32' = '  -1.15%  '
    'D33' = '0.05057'
    'E33' = '  -0.13%  '
    'D34' = '0.7114'
    'E34' = '  -0.61%  '
    'D35' = '1.114'
    'E35' = '  -2.62%  '
    'D36' = '2.669'
    'E36' = '  -0.77%  '
    'D37' = '0.01787'
    'E37' = '  +5.64%  '
    'E38' = '  -4.45%  '
    'D39' = '2.173'
    'E39' = '  -4.30%  '
    'D40' = '0.9298'
    'E40' = '  -3.89%  '
    'B41' = 'TheSandbox'
    'C41' = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    'D41' = '0.4238'
    'E41' = '  -1.34%  '
    'B42' = 'FraxShare'
    'C42' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D42' = '5.760'
    'E42' = '  -6.50%  '
    'B43' = 'PaxDollar'
    'C43' = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
    'D43' = '0.9985'
    'E43' = '  -0.06%  '
    'D44' = '102.71'
    'E44' = '  -2.53%  '
    'D45' = '7.424'
    'E45' = '  -2.80%  '
    'D46' = '0.1268'
    'E46' = '  -0.70%  '
    'D47' = '0.05658'
    'E47' = '  -1.31%  '
    'D48' = '32.50'
    'E48' = '  -1.84%  '
    'D49' = '0.3775'
    'E49' = '  -0.56%  '
    'D50' = '8.275'
    'E50' = '  -1.67%  '
    'D51' = '55.93'
    'E51' = '  -1.29%  '
}

foreach ($addr in $updates.Keys) {
    $col = $addr -replace "[0-9]+$", ""
    # Columns D (Price) and E (Volume) hold numeric-looking / percent-looking
    # text that must stay literal text (keep leading/trailing zeros, the
    # "." thousands separators, the padded "%" strings) instead of being
    # auto-coerced into numbers by the COM layer.
    if ($col -eq "D" -or $col -eq "E") {
        $ws.Range($addr).NumberFormat = "@"
    }
    $ws.Range($addr).Value = $updates[$addr]
}
